# Adds a new "16. 3. 2021" wave of data to both sheets (new last column),
# and bumps the "aktualizace" date in the two footer/title strings from
# 9. 3. 2021 to 23. 3. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": new column AA (after Z), header + 22 rows of percentages
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Copy the formatting of the previous header cell (Z1) onto the new
# header cell (AA1) so the new column matches the existing header style.
$wsData.Range("Z1").Copy()
$wsData.Range("AA1").PasteSpecial(-4122)
$wsData.Range("AA1").Value = "16. 3. 2021"

$dataValues = @{
    2  = 0.13
    3  = 0.1
    4  = 0.16
    5  = 0.11
    6  = 0.1
    7  = 0.19
    8  = 0.14
    9  = 0.28
    10 = 0.13
    11 = 0.09
    12 = 0.17
    13 = 0.09
    14 = 0.28
    15 = 0.16
    16 = 0.09
    17 = 0.21
    18 = 0.14
    19 = 0.09
    20 = 0.1
    21 = 0.09
    22 = 0.06
    23 = 0.19
}

foreach ($row in $dataValues.Keys) {
    $wsData.Cells.Item($row, 27).Value = $dataValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "pocetR": new column Z (after Y), header + 22 rows of counts
# ---------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("pocetR")

$wsCount.Range("Y1").Copy()
$wsCount.Range("Z1").PasteSpecial(-4122)
$wsCount.Range("Z1").Value = "16. 3. 2021"

$countValues = @{
    2  = 2101
    3  = 1015
    4  = 1086
    5  = 292
    6  = 723
    7  = 338
    8  = 748
    9  = 176
    10 = 325
    11 = 401
    12 = 375
    13 = 824
    14 = 244
    15 = 478
    16 = 1379
    17 = 257
    18 = 816
    19 = 615
    20 = 272
    21 = 345
    22 = 738
    23 = 1018
}

foreach ($row in $countValues.Keys) {
    $wsCount.Cells.Item($row, 26).Value = $countValues[$row]
}

# Row 24 on this sheet carries trailing blank cells across every data
# column (B24:Y24) -- extend that blank run to the new Z24 cell too.
$wsCount.Range("Z24").Value = ""

# ---------------------------------------------------------------------
# Bump the "aktualizace 9. 3. 2021" references to "23. 3. 2021" in both
# sheets' title/footer rows (row 24, column A).
# ---------------------------------------------------------------------
$wsData.Range("A24").Value = $wsData.Range("A24").Value().Replace("9. 3. 2021", "23. 3. 2021")
$wsCount.Range("A24").Value = $wsCount.Range("A24").Value().Replace("9. 3. 2021", "23. 3. 2021")
